# Applies the "it works for xls" edit:
#  - workbook.xml: drop activeTab="6" (guess_me, the first sheet, becomes active)
#  - sheet1 ("guess_me"): shift all data one column to the left (B:F -> A:E);
#    the bestFit column width definition moves from column D to column C
#    along with the data. The sheet becomes the tabSelected one and its
#    selection becomes the whole of column A.
#  - sheet7 ("text_coercion"): no longer the tabSelected sheet; its selection
#    moves from A9 to A10.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "guess_me" sheet: delete (now-empty) column A, which shifts B:F left onto
# A:E, keeping formulas, styles and the bestFit column width intact.
# ---------------------------------------------------------------------------
$guessMe = $wb.Worksheets.Item("guess_me")
$guessMe.Activate()
$guessMe.Columns("A").Delete() | Out-Null

# Selection becomes the whole of column A, with A1 as the active cell.
$guessMe.Range("A1:A1048576").Select() | Out-Null

# ---------------------------------------------------------------------------
# "text_coercion" sheet: selection moves from A9 to A10; it is no longer the
# tabSelected sheet (ensured by activating guess_me, above and below).
# ---------------------------------------------------------------------------
$textCoercion = $wb.Worksheets.Item("text_coercion")
$textCoercion.Activate()
$textCoercion.Range("A10").Select() | Out-Null

# Re-activate guess_me last so it ends up as the active/tabSelected sheet.
$guessMe.Activate()
